$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts B:F left to A:E
$ws.Columns.Item(1).Delete()

# Fix the header text for the MODEL_CONDITION -> MODELCONDITION column (now D1)
$ws.Range("D1").Value = "MODELCONDITION"
